# Add 2022-Q3 data
# 1) Update "总计" summary sheet: insert a new row for 2022-Q3 at the top of
#    the data block (row 2), shifting the existing quarters down, and
#    renumber the running index in column A.
# 2) Add a brand-new worksheet named "2022-Q3" right after "总计" holding the
#    per-fund breakdown for the quarter (mirrors the layout used by the other
#    quarterly sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet — insert the 2022-Q3 summary row
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Rows.Item(2).Insert()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 11
$totals.Range("D2").Value = 0.06

$totals.Range("A2").Font.Bold = $true
$totals.Range("A2").Borders.LineStyle = 1
$totals.Range("A2").HorizontalAlignment = -4108
$totals.Range("A2").VerticalAlignment = -4160

# Renumber the running index (column A) for the rows that followed -
# they all shift down by one quarter.
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4
$totals.Range("A7").Value = 5

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet with the fund-level detail
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $totals)
$q3.Name = "2022-Q3"

# Header row
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Fund rows: code, name, scale, stock-position, position-ratio, mkt-value, rank
$rows = @(
    @("257040", "国联安红利混合",              "0.65", "76.59", "5.68", "0.0369", 3),
    @("014861", "申万菱信双禧混合A",            "1.51", "30.74", "0.52", "0.0079", 5),
    @("003717", "中银量化精选灵活配置混合A",    "0.43", "90.72", "1.28", "0.0055", 10),
    @("015466", "太平中证1000指数增强A",        "0.37", "92.23", "0.99", "0.0037", 9),
    @("004988", "人保双利优选混合A",            "0.56", "25.48", "0.55", "0.0031", 6),
    @("013491", "同泰金融精选股票C",            "0.07", "91.55", "3.11", "0.0022", 9),
    @("013490", "同泰金融精选股票A",            "0.02", "91.55", "3.11", "0.0006", 9),
    @("004989", "人保双利优选混合C",            "0.04", "25.48", "0.55", "0.0002", 6),
    @("015467", "太平中证1000指数增强C",        "0.02", "92.23", "0.99", "0.0002", 9),
    @("014862", "申万菱信双禧混合C",            "0.04", "30.74", "0.52", "0.0002", 5),
    @("010484", "中银量化精选灵活配置混合C",    "0.01", "90.72", "1.28", "0.0001", 10)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $r + 2
    $data = $rows[$r]

    $aCell = $q3.Cells.Item($row, 1)
    $aCell.Value = $r
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $bCell = $q3.Cells.Item($row, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $data[0]

    $q3.Cells.Item($row, 3).Value = $data[1]

    $dCell = $q3.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $data[2]

    $eCell = $q3.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $data[3]

    $fCell = $q3.Cells.Item($row, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $data[4]

    $gCell = $q3.Cells.Item($row, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $data[5]

    $q3.Cells.Item($row, 8).Value = $data[6]
}

Write-Host "2022-Q3 sheet added and 总计 updated"
